# Rebuild the "Ementa" table (A1:C24) to match the restructured syllabus layout.
# Strategy: keep every cell that is already correct untouched (so its style /
# shared-string slot survives as-is); only touch the cells the diff actually
# changes, and use Copy+PasteSpecial to pick up column formatting / reuse the
# "01/01/2012" text verbatim (avoids Excel auto-coercing it to a date serial).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose content is removed entirely in the target layout ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()
$ws.Range("B25").Clear()
$ws.Range("C25").Clear()
$ws.Range("B26").Clear()
$ws.Range("C26").Clear()

# --- "01/01/2012" must stay a shared text string, not become a date serial;
#     copy it (and its row-13 B/C formatting) straight from B8/C8 ---
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# --- Existing cells: just overwrite their text, style is left untouched ---
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Range("C23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"

# --- Column A gets the right bold style automatically from the sheet's
#     column defaults, but B/C need their formatting copied explicitly
#     (the overlapping <col> ranges otherwise default new cells to style 1) ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B18").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C18").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"

# --- Row heights: AutoFit first to drop any stale ht/customHeight override,
#     then re-apply the explicit height the target actually wants ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).AutoFit()
$ws.Rows.Item(26).AutoFit()
